$d = $word.ActiveDocument

# Move to the very end of the document body (end of "ha ha ha" paragraph)
$end = $d.Content
$end.Collapse(0)   # wdCollapseEnd

# Start a new paragraph after the last one
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.Move(4, 1)    # wdParagraph, move forward into the new paragraph

# Insert the three runs of text
$end.InsertAfter("Commit")
$end.Collapse(0)
$end.InsertAfter(" by Duc My")
$end.Collapse(0)
$end.InsertAfter(" Nguyen")
$end.Collapse(0)

# Add the _GoBack bookmark at the (now empty) end of the run
$d.Bookmarks.Add("_GoBack", $end)
